$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 401.93995414158786
$ws.Range("C5").Value = 52926.146686182503
$ws.Range("E5").Value = 24122.958278160062
$ws.Range("F5").Value = 757.20375649635287
$ws.Range("G5").Value = 5034.0136633455586
$ws.Range("J5").Value = 83242.262338326065
$ws.Range("B6").Value = 0.000033916809974204643
$ws.Range("C6").Value = 6090.2525177746411
$ws.Range("D6").Value = 3176.6055366253076
$ws.Range("E6").Value = 14063.185452902055
$ws.Range("F6").Value = 241.11632141442743
$ws.Range("H6").Value = 99.635167669819978
$ws.Range("J6").Value = 23670.795030303067
$ws.Range("C7").Value = 6090.2525177746411
$ws.Range("D7").Value = 2446.5685462545644
$ws.Range("E7").Value = 14051.762682717188
$ws.Range("D8").Value = 730.03699037074352
$ws.Range("E8").Value = 11.422770184866787
$ws.Range("H8").Value = 99.635167669819978
$ws.Range("B9").Value = -275.65451944274798
$ws.Range("C9").Value = -47955.704136857777
$ws.Range("D9").Value = -5979.4983924660055
$ws.Range("E9").Value = -28722.282377307893
$ws.Range("F9").Value = -16.021582286988252
$ws.Range("H9").Value = -1261.8843508168598
$ws.Range("J9").Value = -84211.045359178286
$ws.Range("C10").Value = -22787.404014714884
$ws.Range("D10").Value = -5596.7122406302451
$ws.Range("E10").Value = -1679.1112652649533
$ws.Range("C11").Value = -25168.300122142897
$ws.Range("D11").Value = -382.78615183576028
$ws.Range("E11").Value = -27043.171112042939
$ws.Range("H11").Value = -1261.8843508168598
$ws.Range("D12").Value = -1002.590149267652
$ws.Range("J12").Value = -1002.590149267652
$ws.Range("C13").Value = 44.849882055171939
$ws.Range("D13").Value = 44.719879280166758
$ws.Range("J13").Value = 89.569761335338697
$ws.Range("B14").Value = 126.28546861564988
$ws.Range("C14").Value = 11105.544949154541
$ws.Range("D14").Value = -3760.7631258281831
$ws.Range("E14").Value = 9463.861353754226
$ws.Range("F14").Value = 981.73236197407516
$ws.Range("G14").Value = 5034.0136633455586
$ws.Range("H14").Value = -1162.2491831470397
$ws.Range("J14").Value = 21788.425487868826
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = -1.1240450149512071
$ws.Range("E15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = -1.1240450149512071
$ws.Range("B16").Value = 126.28546861564988
$ws.Range("C16").Value = 11105.544949154537
$ws.Range("D16").Value = -3761.8871708431343
$ws.Range("E16").Value = 9463.8613537542296
$ws.Range("F16").Value = 981.73236197407516
$ws.Range("G16").Value = 5034.0136633455586
$ws.Range("H16").Value = -1162.2491831470397
$ws.Range("J16").Value = 21787.301442853874
$ws.Range("C17").Value = -1930.4810265331166
$ws.Range("D17").Value = 1930.4810265331166
$ws.Range("G17").Value = -2092.7724080189364
$ws.Range("H17").Value = 2092.7724080189364
$ws.Range("C18").Value = -9175.0639226214207
$ws.Range("D18").Value = 8940.7948701448513
$ws.Range("E18").Value = -1478.3858772496594
$ws.Range("F18").Value = -544.66443420107191
$ws.Range("G18").Value = -2941.2412553266217
$ws.Range("H18").Value = 2054.8110221894326
$ws.Range("I18").Value = 118.85195238779175
$ws.Range("J18").Value = -3024.8976446766983
$ws.Range("D19").Value = -122.20203054151465
$ws.Range("E19").Value = -1296.428062330529
$ws.Range("F19").Value = -533.195107588528
$ws.Range("G19").Value = -2941.2412553266217
$ws.Range("H19").Value = 2054.8110221894326
$ws.Range("J19").Value = -2838.2554335977607
$ws.Range("C20").Value = -9175.0639226214207
$ws.Range("D20").Value = 9125.3215529965692
$ws.Range("J20").Value = -49.742369624851563
$ws.Range("D21").Value = -62.324652310203362
$ws.Range("E21").Value = -181.95781491913053
$ws.Range("F21").Value = -11.46932661254391
$ws.Range("I21").Value = 118.85195238779175
$ws.Range("J21").Value = -136.89984145408604
$ws.Range("D22").Value = 581.86315444322327
$ws.Range("E22").Value = 3761.5887271162856
$ws.Range("H22").Value = 430.6151697031857
$ws.Range("I22").Value = 30.554445938952057
$ws.Range("J22").Value = 4804.6214972016469
$ws.Range("B23").Value = 126.28546861564988
$ws.Range("D23").Value = 6527.5255713916094
$ws.Range("E23").Value = 4223.8867493882835
$ws.Range("F23").Value = 437.06792777300325
$ws.Range("H23").Value = 2554.7190773581437
$ws.Range("I23").Value = 88.297506448839698
$ws.Range("J23").Value = 13957.78230097553
$ws.Range("D24").Value = 1125.2647513822537
$ws.Range("E24").Value = 59.263883971473689
$ws.Range("I24").Value = 52.619899118795843
$ws.Range("J24").Value = 1237.1485344725234
$ws.Range("B25").Value = 80.37196056133088
$ws.Range("D25").Value = 177.56837512603184
$ws.Range("E25").Value = 935.55392882101103
$ws.Range("F25").Value = 98.132461902953565
$ws.Range("H25").Value = 714.94006808840356
$ws.Range("I25").Value = 9.4953201417375954
$ws.Range("J25").Value = 2016.0621146414685
$ws.Range("B26").Value = 43.430782459157591
$ws.Range("D26").Value = 263.14447083901541
$ws.Range("E26").Value = 2355.44966903881
$ws.Range("F26").Value = 160.10730182375826
$ws.Range("H26").Value = 1059.1373970396792
$ws.Range("I26").Value = 26.182287188306248
$ws.Range("J26").Value = 3907.4519083887267
$ws.Range("B27").Value = 0.36113499570077584
$ws.Range("D27").Value = 4449.0157472838027
$ws.Range("F27").Value = 92.185039969441235
$ws.Range("H27").Value = 38.348248213293964
$ws.Range("J27").Value = 4579.910170462239
$ws.Range("B28").Value = 2.1129921471820268
$ws.Range("D28").Value = 512.53222676050734
$ws.Range("E28").Value = 873.61926755698801
$ws.Range("F28").Value = 86.643124076850185
$ws.Range("H28").Value = 742.29336401676699
$ws.Range("J28").Value = 2217.2009745582945

$ws.Range("B5:J28").Select()
